$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-11-27 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-11-28 Thursday", 2)
$d.Content.Find.Execute("752÷8=94, 0", $true, $false, $false, $false, $false, $true, 1, $false, "476÷2=238, 0", 2)
$d.Content.Find.Execute("201÷4=50, 1", $true, $false, $false, $false, $false, $true, 1, $false, "309÷7=44, 1", 2)
$d.Content.Find.Execute("844÷8=105, 4", $true, $false, $false, $false, $false, $true, 1, $false, "359÷3=119, 2", 2)
$d.Content.Find.Execute("338÷8=42, 2", $true, $false, $false, $false, $false, $true, 1, $false, "844÷4=211, 0", 2)
$d.Content.Find.Execute("401÷2=200, 1", $true, $false, $false, $false, $false, $true, 1, $false, "699÷5=139, 4", 2)
$d.Content.Find.Execute("444÷8=55, 4", $true, $false, $false, $false, $false, $true, 1, $false, "894÷4=223, 2", 2)
$d.Content.Find.Execute("880÷4=220, 0", $true, $false, $false, $false, $false, $true, 1, $false, "241÷4=60, 1", 2)
$d.Content.Find.Execute("135÷5=27, 0", $true, $false, $false, $false, $false, $true, 1, $false, "124÷4=31, 0", 2)
$d.Content.Find.Execute("675÷9=75, 0", $true, $false, $false, $false, $false, $true, 1, $false, "181÷7=25, 6", 2)
$d.Content.Find.Execute("680÷6=113, 2", $true, $false, $false, $false, $false, $true, 1, $false, "397÷2=198, 1", 2)
$d.Content.Find.Execute("292÷6=48, 4", $true, $false, $false, $false, $false, $true, 1, $false, "234÷4=58, 2", 2)
$d.Content.Find.Execute("468÷4=117, 0", $true, $false, $false, $false, $false, $true, 1, $false, "196÷4=49, 0", 2)
$d.Content.Find.Execute("325÷3=108, 1", $true, $false, $false, $false, $false, $true, 1, $false, "100÷6=16, 4", 2)
$d.Content.Find.Execute("188÷4=47, 0", $true, $false, $false, $false, $false, $true, 1, $false, "302÷2=151, 0", 2)
$d.Content.Find.Execute("556÷8=69, 4", $true, $false, $false, $false, $false, $true, 1, $false, "791÷2=395, 1", 2)
$d.Content.Find.Execute("853÷5=170, 3", $true, $false, $false, $false, $false, $true, 1, $false, "103÷9=11, 4", 2)
$d.Content.Find.Execute("847÷4=211, 3", $true, $false, $false, $false, $false, $true, 1, $false, "802÷9=89, 1", 2)
$d.Content.Find.Execute("546÷8=68, 2", $true, $false, $false, $false, $false, $true, 1, $false, "480÷7=68, 4", 2)
$d.Content.Find.Execute("113÷8=14, 1", $true, $false, $false, $false, $false, $true, 1, $false, "629÷8=78, 5", 2)
$d.Content.Find.Execute("587÷2=293, 1", $true, $false, $false, $false, $false, $true, 1, $false, "957÷2=478, 1", 2)
$d.Content.Find.Execute("386÷6=64, 2", $true, $false, $false, $false, $false, $true, 1, $false, "678÷3=226, 0", 2)
$d.Content.Find.Execute("324÷9=36, 0", $true, $false, $false, $false, $false, $true, 1, $false, "615÷7=87, 6", 2)
$d.Content.Find.Execute("270÷4=67, 2", $true, $false, $false, $false, $false, $true, 1, $false, "328÷5=65, 3", 2)
$d.Content.Find.Execute("107÷7=15, 2", $true, $false, $false, $false, $false, $true, 1, $false, "272÷5=54, 2", 2)
$d.Content.Find.Execute("487÷2=243, 1", $true, $false, $false, $false, $false, $true, 1, $false, "390÷2=195, 0", 2)
